$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# Overview sheet: E2 and F2 both use the "Ready for handoff" string.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: C2 uses the "Ready for handoff" string.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "In Translation"

# de-de sheet: C2 uses the "Ready for handoff" string.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "In Translation"

# --- Column width changes ---
# Overview sheet: columns E and F narrower.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C narrower.
$wsZh.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C narrower.
$wsDe.Columns.Item(3).ColumnWidth = 12.5
